$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = [double]"0.01253208636536152"
$ws.Range("C2").Value = [double]"6.708468553440206e-05"
$ws.Range("D2").Value = [double]"0.1496068669990043"
$ws.Range("E2").Value = [double]"13.86384647080068"
$ws.Range("G2").Value = [double]"14.02605250885058"

# Row 3
$ws.Range("B3").Value = [double]"3.272327238179451"
$ws.Range("C3").Value = [double]"1.626987699542094"
$ws.Range("D3").Value = [double]"0.1496068669990043"
$ws.Range("E3").Value = [double]"0.5333859586016987"
$ws.Range("G3").Value = [double]"5.582307763322248"

# Row 4
$ws.Range("B4").Value = [double]"3.272327238179451"
$ws.Range("C4").Value = [double]"1.626987699542094"
$ws.Range("D4").Value = [double]"3.223369029078222"
$ws.Range("E4").Value = [double]"0.5333859586016987"
$ws.Range("G4").Value = [double]"8.656069925401464"
